# Update NATMI ligand/receptor expression + specificity + edge-weight
# columns (G:T, rows 2:26) with recomputed values from the new TPM input
# (ligand = Apoe, receptor = Ldlr). Columns A:F (cluster/gene metadata,
# ligand-expressing cells, detection rate) are unchanged by this edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 25,14
$arr[0,0] = 32.50235
$arr[0,1] = 97.50704999999999
$arr[0,2] = 0.004318312013857221
$arr[0,3] = 0.004318312013857221
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 2.590549
$arr[0,7] = 7.771647
$arr[0,8] = 0.08453422544559429
$arr[0,9] = 0.0845342254455943
$arr[0,10] = 84.19893029014999
$arr[0,11] = 757.7903726113499
$arr[0,12] = 0.0003650451613238246
$arr[0,13] = 0.0003650451613238247
$arr[1,0] = 32.50235
$arr[1,1] = 97.50704999999999
$arr[1,2] = 0.004318312013857221
$arr[1,3] = 0.004318312013857221
$arr[1,4] = 3
$arr[1,5] = 1
$arr[1,6] = 6.056162
$arr[1,7] = 18.168486
$arr[1,8] = 0.1976233469596758
$arr[1,9] = 0.1976233469596758
$arr[1,10] = 196.8394969807
$arr[1,11] = 1771.5554728263
$arr[1,12] = 0.0008533992733946421
$arr[1,13] = 0.0008533992733946421
$arr[2,0] = 32.50235
$arr[2,1] = 97.50704999999999
$arr[2,2] = 0.004318312013857221
$arr[2,3] = 0.004318312013857221
$arr[2,4] = 3
$arr[2,5] = 1
$arr[2,6] = 12.64302866666667
$arr[2,7] = 37.929086
$arr[2,8] = 0.4125645319286033
$arr[2,9] = 0.4125645319286034
$arr[2,10] = 410.9281427840333
$arr[2,11] = 3698.3532850563
$arr[2,12] = 0.001781582374718669
$arr[2,13] = 0.001781582374718669
$arr[3,0] = 32.50235
$arr[3,1] = 97.50704999999999
$arr[3,2] = 0.004318312013857221
$arr[3,3] = 0.004318312013857221
$arr[3,4] = 3
$arr[3,5] = 1
$arr[3,6] = 3.370263
$arr[3,7] = 10.110789
$arr[3,8] = 0.1099776812764186
$arr[3,9] = 0.1099776812764186
$arr[3,10] = 109.54146761805
$arr[3,11] = 985.8732085624499
$arr[3,12] = 0.0004749179423121189
$arr[3,13] = 0.000474917942312119
$arr[4,0] = 32.50235
$arr[4,1] = 97.50704999999999
$arr[4,2] = 0.004318312013857221
$arr[4,3] = 0.004318312013857221
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 5.984969666666667
$arr[4,7] = 17.954909
$arr[4,8] = 0.1953002143897079
$arr[4,9] = 0.1953002143897079
$arr[4,10] = 194.5255788453833
$arr[4,11] = 1750.73020960845
$arr[4,12] = 0.0008433672621079664
$arr[4,13] = 0.0008433672621079667
$arr[5,0] = 50.89916233333333
$arr[5,1] = 152.697487
$arr[5,2] = 0.006762540683959845
$arr[5,3] = 0.006762540683959845
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 2.590549
$arr[5,7] = 7.771647
$arr[5,8] = 0.08453422544559429
$arr[5,9] = 0.0845342254455943
$arr[5,10] = 131.8567740834543
$arr[5,11] = 1186.710966751089
$arr[5,12] = 0.000571666138762865
$arr[5,13] = 0.0005716661387628651
$arr[6,0] = 50.89916233333333
$arr[6,1] = 152.697487
$arr[6,2] = 0.006762540683959845
$arr[6,3] = 0.006762540683959845
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 6.056162
$arr[6,7] = 18.168486
$arr[6,8] = 0.1976233469596758
$arr[6,9] = 0.1976233469596758
$arr[6,10] = 308.2535727549647
$arr[6,11] = 2774.282154794682
$arr[6,12] = 0.00133643592391512
$arr[6,13] = 0.00133643592391512
$arr[7,0] = 50.89916233333333
$arr[7,1] = 152.697487
$arr[7,2] = 0.006762540683959845
$arr[7,3] = 0.006762540683959845
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 12.64302866666667
$arr[7,7] = 37.929086
$arr[7,8] = 0.4125645319286033
$arr[7,9] = 0.4125645319286034
$arr[7,10] = 643.5195684896535
$arr[7,11] = 5791.676116406881
$arr[7,12] = 0.002789984431926031
$arr[7,13] = 0.002789984431926031
$arr[8,0] = 50.89916233333333
$arr[8,1] = 152.697487
$arr[8,2] = 0.006762540683959845
$arr[8,3] = 0.006762540683959845
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 3.370263
$arr[8,7] = 10.110789
$arr[8,8] = 0.1099776812764186
$arr[8,9] = 0.1099776812764186
$arr[8,10] = 171.543563543027
$arr[8,11] = 1543.892071887243
$arr[8,12] = 0.0007437285439593498
$arr[8,13] = 0.0007437285439593499
$arr[9,0] = 50.89916233333333
$arr[9,1] = 152.697487
$arr[9,2] = 0.006762540683959845
$arr[9,3] = 0.006762540683959845
$arr[9,4] = 3
$arr[9,5] = 1
$arr[9,6] = 5.984969666666667
$arr[9,7] = 17.954909
$arr[9,8] = 0.1953002143897079
$arr[9,9] = 0.1953002143897079
$arr[9,10] = 304.6299426237426
$arr[9,11] = 2741.669483613683
$arr[9,12] = 0.00132072564539648
$arr[9,13] = 0.00132072564539648
$arr[10,0] = 3274.382486666667
$arr[10,1] = 9823.14746
$arr[10,2] = 0.4350394734576531
$arr[10,3] = 0.435039473457653
$arr[10,4] = 3
$arr[10,5] = 1
$arr[10,6] = 2.590549
$arr[10,7] = 7.771647
$arr[10,8] = 0.08453422544559429
$arr[10,9] = 0.0845342254455943
$arr[10,10] = 8482.448276451847
$arr[10,11] = 76342.03448806662
$arr[10,12] = 0.03677572492700188
$arr[10,13] = 0.03677572492700188
$arr[11,0] = 3274.382486666667
$arr[11,1] = 9823.14746
$arr[11,2] = 0.4350394734576531
$arr[11,3] = 0.435039473457653
$arr[11,4] = 3
$arr[11,5] = 1
$arr[11,6] = 6.056162
$arr[11,7] = 18.168486
$arr[11,8] = 0.1976233469596758
$arr[11,9] = 0.1976233469596758
$arr[11,10] = 19830.19078921618
$arr[11,11] = 178471.7171029456
$arr[11,12] = 0.08597395680427647
$arr[11,13] = 0.08597395680427646
$arr[12,0] = 3274.382486666667
$arr[12,1] = 9823.14746
$arr[12,2] = 0.4350394734576531
$arr[12,3] = 0.435039473457653
$arr[12,4] = 3
$arr[12,5] = 1
$arr[12,6] = 12.64302866666667
$arr[12,7] = 37.929086
$arr[12,8] = 0.4125645319286033
$arr[12,9] = 0.4125645319286034
$arr[12,10] = 41398.11164455795
$arr[12,11] = 372583.0048010215
$arr[12,12] = 0.1794818567375227
$arr[12,13] = 0.1794818567375227
$arr[13,0] = 3274.382486666667
$arr[13,1] = 9823.14746
$arr[13,2] = 0.4350394734576531
$arr[13,3] = 0.435039473457653
$arr[13,4] = 3
$arr[13,5] = 1
$arr[13,6] = 3.370263
$arr[13,7] = 10.110789
$arr[13,8] = 0.1099776812764186
$arr[13,9] = 0.1099776812764186
$arr[13,10] = 11035.53014266066
$arr[13,11] = 99319.77128394594
$arr[13,12] = 0.04784463255458675
$arr[13,13] = 0.04784463255458675
$arr[14,0] = 3274.382486666667
$arr[14,1] = 9823.14746
$arr[14,2] = 0.4350394734576531
$arr[14,3] = 0.435039473457653
$arr[14,4] = 3
$arr[14,5] = 1
$arr[14,6] = 5.984969666666667
$arr[14,7] = 17.954909
$arr[14,8] = 0.1953002143897079
$arr[14,9] = 0.1953002143897079
$arr[14,10] = 19597.07985976457
$arr[14,11] = 176373.7187378811
$arr[14,12] = 0.08496330243426528
$arr[14,13] = 0.08496330243426528
$arr[15,0] = 7.278837333333333
$arr[15,1] = 21.836512
$arr[15,2] = 0.0009670774791190726
$arr[15,3] = 0.0009670774791190726
$arr[15,4] = 3
$arr[15,5] = 1
$arr[15,6] = 2.590549
$arr[15,7] = 7.771647
$arr[15,8] = 0.08453422544559429
$arr[15,9] = 0.0845342254455943
$arr[15,10] = 18.85618477502933
$arr[15,11] = 169.705662975264
$arr[15,12] = 0.00008175114564320869
$arr[15,13] = 0.0000817511456432087
$arr[16,0] = 7.278837333333333
$arr[16,1] = 21.836512
$arr[16,2] = 0.0009670774791190726
$arr[16,3] = 0.0009670774791190726
$arr[16,4] = 3
$arr[16,5] = 1
$arr[16,6] = 6.056162
$arr[16,7] = 18.168486
$arr[16,8] = 0.1976233469596758
$arr[16,9] = 0.1976233469596758
$arr[16,10] = 44.08181806231467
$arr[16,11] = 396.736362560832
$arr[16,12] = 0.0001911170881928372
$arr[16,13] = 0.0001911170881928372
$arr[17,0] = 7.278837333333333
$arr[17,1] = 21.836512
$arr[17,2] = 0.0009670774791190726
$arr[17,3] = 0.0009670774791190726
$arr[17,4] = 3
$arr[17,5] = 1
$arr[17,6] = 12.64302866666667
$arr[17,7] = 37.929086
$arr[17,8] = 0.4125645319286033
$arr[17,9] = 0.4125645319286034
$arr[17,10] = 92.02654906533688
$arr[17,11] = 828.2389415880319
$arr[17,12] = 0.0003989818675114539
$arr[17,13] = 0.0003989818675114539
$arr[18,0] = 7.278837333333333
$arr[18,1] = 21.836512
$arr[18,2] = 0.0009670774791190726
$arr[18,3] = 0.0009670774791190726
$arr[18,4] = 3
$arr[18,5] = 1
$arr[18,6] = 3.370263
$arr[18,7] = 10.110789
$arr[18,8] = 0.1099776812764186
$arr[18,9] = 0.1099776812764186
$arr[18,10] = 24.531596147552
$arr[18,11] = 220.784365327968
$arr[18,12] = 0.0001063569387681597
$arr[18,13] = 0.0001063569387681598
$arr[19,0] = 7.278837333333333
$arr[19,1] = 21.836512
$arr[19,2] = 0.0009670774791190726
$arr[19,3] = 0.0009670774791190726
$arr[19,4] = 3
$arr[19,5] = 1
$arr[19,6] = 5.984969666666667
$arr[19,7] = 17.954909
$arr[19,8] = 0.1953002143897079
$arr[19,9] = 0.1953002143897079
$arr[19,10] = 43.56362064860089
$arr[19,11] = 392.072585837408
$arr[19,12] = 0.0001888704390034131
$arr[19,13] = 0.0001888704390034132
$arr[20,0] = 4161.570231333333
$arr[20,1] = 12484.710694
$arr[20,2] = 0.5529125963654108
$arr[20,3] = 0.5529125963654108
$arr[20,4] = 3
$arr[20,5] = 1
$arr[20,6] = 2.590549
$arr[20,7] = 7.771647
$arr[20,8] = 0.08453422544559429
$arr[20,9] = 0.0845342254455943
$arr[20,10] = 10780.75160121033
$arr[20,11] = 97026.76441089301
$arr[20,12] = 0.04674003807286251
$arr[20,13] = 0.04674003807286252
$arr[21,0] = 4161.570231333333
$arr[21,1] = 12484.710694
$arr[21,2] = 0.5529125963654108
$arr[21,3] = 0.5529125963654108
$arr[21,4] = 3
$arr[21,5] = 1
$arr[21,6] = 6.056162
$arr[21,7] = 18.168486
$arr[21,8] = 0.1976233469596758
$arr[21,9] = 0.1976233469596758
$arr[21,10] = 25203.14349533214
$arr[21,11] = 226828.2914579893
$arr[21,12] = 0.1092684378698968
$arr[21,13] = 0.1092684378698968
$arr[22,0] = 4161.570231333333
$arr[22,1] = 12484.710694
$arr[22,2] = 0.5529125963654108
$arr[22,3] = 0.5529125963654108
$arr[22,4] = 3
$arr[22,5] = 1
$arr[22,6] = 12.64302866666667
$arr[22,7] = 37.929086
$arr[22,8] = 0.4125645319286033
$arr[22,9] = 0.4125645319286034
$arr[22,10] = 52614.85173309396
$arr[22,11] = 473533.6655978456
$arr[22,12] = 0.2281121265169245
$arr[22,13] = 0.2281121265169245
$arr[23,0] = 4161.570231333333
$arr[23,1] = 12484.710694
$arr[23,2] = 0.5529125963654108
$arr[23,3] = 0.5529125963654108
$arr[23,4] = 3
$arr[23,5] = 1
$arr[23,6] = 3.370263
$arr[23,7] = 10.110789
$arr[23,8] = 0.1099776812764186
$arr[23,9] = 0.1099776812764186
$arr[23,10] = 14025.58617256417
$arr[23,11] = 126230.2755530776
$arr[23,12] = 0.06080804529679224
$arr[23,13] = 0.06080804529679225
$arr[24,0] = 4161.570231333333
$arr[24,1] = 12484.710694
$arr[24,2] = 0.5529125963654108
$arr[24,3] = 0.5529125963654108
$arr[24,4] = 3
$arr[24,5] = 1
$arr[24,6] = 5.984969666666667
$arr[24,7] = 17.954909
$arr[24,8] = 0.1953002143897079
$arr[24,9] = 0.1953002143897079
$arr[24,10] = 24906.87160023299
$arr[24,11] = 224161.8444020968
$arr[24,12] = 0.1079839486089347
$arr[24,13] = 0.1079839486089348
$ws.Range("G2:T26").Value = $arr
Write-Output "applied new TPM values to G2:T26"
